# The weekly refresh reshuffles the date/variety/quality/price rows for the
# Damasco - Agro Chillan subset. Columns A,B,C and E-J are constant across
# every row already, so only D and K:T need to move. We snapshot the
# current values for D and K:T on rows 2-12, then write them back out in
# the new (permuted) row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 12

# Columns that actually vary row to row in this subset.
$cols = @("D","K","L","M","N","O","P","Q","R","S","T")

# Snapshot current ("before") values per row/column.
# NB: use .Value2 (not .Value) to read - this host's COM shim mishandles
# chained/stored reads through the .Value getter.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Maps each SOURCE (before) row to its DESTINATION (after) row.
$rowMap = @{
    2  = 7
    3  = 9
    4  = 10
    5  = 4
    6  = 2
    7  = 3
    8  = 5
    9  = 12
    10 = 6
    11 = 8
    12 = 11
}

foreach ($srcRow in $rowMap.Keys) {
    $dstRow = $rowMap[$srcRow]
    $rowVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$dstRow").Value2 = $rowVals[$c]
    }
}
